$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1. "The new version of AutoTracker run all the same calculations..."
#    -> split into three runs: "...run" + "s" + " all the same..."
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(3).Range
$xml1 = "<w:p $wns>" +
        "<w:r><w:t>The new version of AutoTracker run</w:t></w:r>" +
        "<w:r><w:t>s</w:t></w:r>" +
        '<w:r><w:t xml:space="preserve"> all the same calculations as the previous version while introducing refinements that streamline use, improve processing speed, and minimize errors.</w:t></w:r>' +
        "</w:p>"
$p1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2. "The new AutoTracker code is now formatted as Python module..."
#    -> rewritten / split into multiple runs, with a spell-checked "hilllab"
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(5).Range
$xml2 = "<w:p $wns>" +
        '<w:r><w:t xml:space="preserve">The new AutoTracker code is now </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">a part of the </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>hilllab</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> Python module. </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">To import the AutoTracker, run the code below. </w:t></w:r>' +
        "</w:p>"
$p2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3. "File paths are chosen using the function below. It opens a dialog to
#    select a directory..." -> split first run into three runs, keep the
#    remaining two runs of the paragraph untouched.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(9).Range
$xml3 = "<w:p $wns>" +
        '<w:r><w:t>File paths are chosen using the function below. It opens a dialog</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:t>to select a directory containing videos and another dialog to choose where to save the VRPN files.</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:t>It also prints a summary so you can verify your selections.</w:t></w:r>' +
        "</w:p>"
$p3.InsertXML($xml3)

# ---------------------------------------------------------------------------
# 4. Table cell "autotrack_videos_parameter_test()" -> merge the three runs
#    into a single run (keep Quote style / color=auto formatting).
# ---------------------------------------------------------------------------
$cellPara = $d.Tables(3).Cell(1, 1).Range.Paragraphs(1).Range
$xml4 = "<w:p $wns>" +
        '<w:pPr><w:pStyle w:val="Quote"/><w:rPr><w:color w:val="auto"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>autotrack_videos_parameter_test()</w:t></w:r>' +
        "</w:p>"
$cellPara.InsertXML($xml4)

# ---------------------------------------------------------------------------
# 5. Remove the trailing page break plus the entire "VRPN Structure" section
#    (heading, two paragraphs, and the column-reference table) so that the
#    "Once running..." paragraph becomes the last paragraph of the document.
# ---------------------------------------------------------------------------
# Drop the VRPN-column-reference table entirely.
$d.Tables($d.Tables.Count).Delete()

# Locate the "Once running..." paragraph and the page-break run right after it.
$onceRange = $d.Content
$onceRange.Find.Execute("Once running, the AutoTracker saves VRPN files to the selected location. If it encounters certain errors while processing a video, it will skip that video and create a report in the video" + [char]8217 + "s folder describing the issue.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterText = $onceRange.End

# Delete the page-break run immediately following the matched sentence. After
# this delete, the paragraph mark that used to sit right after the page break
# shifts back to $afterText, so it is still there (and must be preserved).
$d.Range($afterText, $afterText + 1).Delete()

# Delete everything from just after that (now-preserved) paragraph mark through
# the end of the document (the "VRPN Structure" heading/paragraphs, the table,
# and the trailing empty paragraph).
$paraEnd = $afterText + 1
$docEnd = $d.Content.End
if ($docEnd -gt $paraEnd) {
    $d.Range($paraEnd, $docEnd).Delete()
}
